$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing Table1 occupies A1:M2. Insert a new column before it (shifts
# the table right to B1:N2) so a new "Season" column can be placed in A.
$tbl = $ws.ListObjects.Item(1)
$tbl.Range.Columns(1).EntireColumn.Insert()
$tbl.Resize($ws.Range("B1:N2"))

# New "Season" column values for the existing row.
$ws.Range("A1").Value = "Season"
$ws.Range("A2").Value = "24/25"

# Add a second row of stats (23/24 season) to the table.
$newRow = $tbl.ListRows.Add()
$ws.Range("A3").Value = "23/24"
$ws.Range("B3").Value = "Arsenal"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 89
$ws.Range("E3").Value = 38
$ws.Range("F3").Value = 28
$ws.Range("G3").Value = 5
$ws.Range("H3").Value = 5
$ws.Range("I3").Value = 91
$ws.Range("J3").Value = 29
$ws.Range("K3").Value = 62
$ws.Range("L3").Value = "Bukayo Saka"
$ws.Range("M3").Value = 2
$ws.Range("N3").Value = 62

# Match the "Top Scorer" cell formatting used in row 2.
$ws.Range("L2").Copy()
$ws.Range("L3").PasteSpecial(-4122)
$ws.Range("L3").Value = "Bukayo Saka"

# Apply a filter on the new Season column (outside of Table1's own filter).
$ws.Range("A1:A3").AutoFilter()

$nm = $ws.Names.Add("_xlnm._FilterDatabase", "='Arsenal Stats'!`$A`$1:`$A`$3")
$nm.Visible = $false

$ws.Range("C13").Select()
